# Add list of mac addresses of demo boards
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Boards")

# Each entry: RowNumber -> (Assembly, Chip ID (hex), Mac Address)
# Column D (Flash ID (hex)) is constant "0x4016" for every board.
$rows = @(
    @{ Row = 2;  Assembly = "Lite"; Chip = "0x00d57af7"; Mac = "18:fe:34:d5:7a:f7" },
    @{ Row = 4;  Assembly = "Lite"; Chip = "0x00d578d2"; Mac = "18:fe:34:d5:78:d2" },
    @{ Row = 5;  Assembly = "Lite"; Chip = "0x00d767c7"; Mac = "18:fe:34:d7:67:c7" },
    @{ Row = 6;  Assembly = "Lite"; Chip = "0x001487d5"; Mac = "5c:cf:7f:14:87:d5" },
    @{ Row = 8;  Assembly = "Lite"; Chip = "0x00d76830"; Mac = "18:fe:34:d7:68:30" },
    @{ Row = 9;  Assembly = "Lite"; Chip = "0x00d57b86"; Mac = "18:fe:34:d5:7b:86" },
    @{ Row = 11; Assembly = "Lite"; Chip = "0x00d76879"; Mac = "18:fe:34:d7:68:79" },
    @{ Row = 12; Assembly = "Lite"; Chip = "0x00149b1f"; Mac = "5c:cf:7f:14:9b:1f" },
    @{ Row = 13; Assembly = "Lite"; Chip = "0x000a7895"; Mac = " 5c:cf:7f:0a:78:95" }
)

# A row that already carries the "data row" number/mono-space formatting in
# columns C:E, used as the formatting source for rows whose C:E cells don't
# exist yet (rows 2 and 4 started completely empty beyond column A).
# Re-pasting the same format onto rows that already have it is a harmless
# no-op (it does not create duplicate style entries).
$formatSource = $ws.Range("C5:E5")

foreach ($r in $rows) {
    $rowNum = $r.Row
    $destRange = $ws.Range("C" + $rowNum + ":E" + $rowNum)
    $formatSource.Copy()
    $destRange.PasteSpecial(-4122)

    $ws.Cells.Item($rowNum, 2).Value = $r.Assembly
    $ws.Cells.Item($rowNum, 3).Value = $r.Chip
    $ws.Cells.Item($rowNum, 4).Value = "0x4016"
    $ws.Cells.Item($rowNum, 5).Value = $r.Mac
}

# Rows 14 and 15 only pick up the "Full" assembly marker; C:E stay blank.
$ws.Cells.Item(14, 2).Value = "Full"
$ws.Cells.Item(15, 2).Value = "Full"

$excel.CutCopyMode = $false

# Restore selection to E2 (matches the final state captured in the diff).
$ws.Range("E2").Select()
